$d = $word.ActiveDocument
$d.Content.Find.Execute("90-88=", $true, $true, $false, $false, $false, $true, 1, $false, "13+81=", 2) | Out-Null
$d.Content.Find.Execute("99-78=", $true, $true, $false, $false, $false, $true, 1, $false, "30+11=", 2) | Out-Null
$d.Content.Find.Execute("39-24=", $true, $true, $false, $false, $false, $true, 1, $false, "17+19=", 2) | Out-Null
$d.Content.Find.Execute("72-51=", $true, $true, $false, $false, $false, $true, 1, $false, "84-40=", 2) | Out-Null
$d.Content.Find.Execute("66-63=", $true, $true, $false, $false, $false, $true, 1, $false, "76-61=", 2) | Out-Null
$d.Content.Find.Execute("46+53=", $true, $true, $false, $false, $false, $true, 1, $false, "5+0=", 2) | Out-Null
$d.Content.Find.Execute("61-12=", $true, $true, $false, $false, $false, $true, 1, $false, "86-35=", 2) | Out-Null
$d.Content.Find.Execute("16+20=", $true, $true, $false, $false, $false, $true, 1, $false, "31-1=", 2) | Out-Null
$d.Content.Find.Execute("44-8=", $true, $true, $false, $false, $false, $true, 1, $false, "38-5=", 2) | Out-Null
$d.Content.Find.Execute("67+5=", $true, $true, $false, $false, $false, $true, 1, $false, "48+14=", 2) | Out-Null
$d.Content.Find.Execute("11+50=", $true, $true, $false, $false, $false, $true, 1, $false, "67-0=", 2) | Out-Null
$d.Content.Find.Execute("4+81=", $true, $true, $false, $false, $false, $true, 1, $false, "97-96=", 2) | Out-Null
$d.Content.Find.Execute("16-5=", $true, $true, $false, $false, $false, $true, 1, $false, "13+9=", 2) | Out-Null
$d.Content.Find.Execute("62+6=", $true, $true, $false, $false, $false, $true, 1, $false, "2+67=", 2) | Out-Null
$d.Content.Find.Execute("69-20=", $true, $true, $false, $false, $false, $true, 1, $false, "38+51=", 2) | Out-Null
$d.Content.Find.Execute("85-64=", $true, $true, $false, $false, $false, $true, 1, $false, "35+55=", 2) | Out-Null
$d.Content.Find.Execute("92-63=", $true, $true, $false, $false, $false, $true, 1, $false, "64-51=", 2) | Out-Null
$d.Content.Find.Execute("0+56=", $true, $true, $false, $false, $false, $true, 1, $false, "8+22=", 2) | Out-Null
$d.Content.Find.Execute("43-41=", $true, $true, $false, $false, $false, $true, 1, $false, "7+88=", 2) | Out-Null
$d.Content.Find.Execute("84-51=", $true, $true, $false, $false, $false, $true, 1, $false, "76-68=", 2) | Out-Null
$d.Content.Find.Execute("14+78=", $true, $true, $false, $false, $false, $true, 1, $false, "40+9=", 2) | Out-Null
$d.Content.Find.Execute("69-32=", $true, $true, $false, $false, $false, $true, 1, $false, "80-53=", 2) | Out-Null
$d.Content.Find.Execute("46+6=", $true, $true, $false, $false, $false, $true, 1, $false, "61+25=", 2) | Out-Null
$d.Content.Find.Execute("99-65=", $true, $true, $false, $false, $false, $true, 1, $false, "50-0=", 2) | Out-Null
$d.Content.Find.Execute("28+9=", $true, $true, $false, $false, $false, $true, 1, $false, "83-15=", 2) | Out-Null
$d.Content.Find.Execute("60-11=", $true, $true, $false, $false, $false, $true, 1, $false, "72+10=", 2) | Out-Null
$d.Content.Find.Execute("82-37=", $true, $true, $false, $false, $false, $true, 1, $false, "52-12=", 2) | Out-Null
$d.Content.Find.Execute("15+68=", $true, $true, $false, $false, $false, $true, 1, $false, "28-25=", 2) | Out-Null
$d.Content.Find.Execute("12+83=", $true, $true, $false, $false, $false, $true, 1, $false, "41+17=", 2) | Out-Null
$d.Content.Find.Execute("2+60=", $true, $true, $false, $false, $false, $true, 1, $false, "76-39=", 2) | Out-Null
$d.Content.Find.Execute("9+84=", $true, $true, $false, $false, $false, $true, 1, $false, "87-15=", 2) | Out-Null
$d.Content.Find.Execute("81+15=", $true, $true, $false, $false, $false, $true, 1, $false, "75-0=", 2) | Out-Null
$d.Content.Find.Execute("96-37=", $true, $true, $false, $false, $false, $true, 1, $false, "43+38=", 2) | Out-Null
$d.Content.Find.Execute("40+6=", $true, $true, $false, $false, $false, $true, 1, $false, "75-12=", 2) | Out-Null
$d.Content.Find.Execute("47+23=", $true, $true, $false, $false, $false, $true, 1, $false, "3+36=", 2) | Out-Null
$d.Content.Find.Execute("48-1=", $true, $true, $false, $false, $false, $true, 1, $false, "67-33=", 2) | Out-Null
$d.Content.Find.Execute("95-51=", $true, $true, $false, $false, $false, $true, 1, $false, "74-57=", 2) | Out-Null
$d.Content.Find.Execute("10+60=", $true, $true, $false, $false, $false, $true, 1, $false, "52-31=", 2) | Out-Null
$d.Content.Find.Execute("29+40=", $true, $true, $false, $false, $false, $true, 1, $false, "93-44=", 2) | Out-Null
$d.Content.Find.Execute("42+28=", $true, $true, $false, $false, $false, $true, 1, $false, "9+40=", 2) | Out-Null
$d.Content.Find.Execute("17+47=", $true, $true, $false, $false, $false, $true, 1, $false, "94-43=", 2) | Out-Null
$d.Content.Find.Execute("24+24=", $true, $true, $false, $false, $false, $true, 1, $false, "48-37=", 2) | Out-Null
$d.Content.Find.Execute("51+42=", $true, $true, $false, $false, $false, $true, 1, $false, "45+32=", 2) | Out-Null
$d.Content.Find.Execute("38+7=", $true, $true, $false, $false, $false, $true, 1, $false, "54-52=", 2) | Out-Null
$d.Content.Find.Execute("11+12=", $true, $true, $false, $false, $false, $true, 1, $false, "53-23=", 2) | Out-Null
$d.Content.Find.Execute("87-51=", $true, $true, $false, $false, $false, $true, 1, $false, "79-25=", 2) | Out-Null
$d.Content.Find.Execute("62+19=", $true, $true, $false, $false, $false, $true, 1, $false, "24+72=", 2) | Out-Null
$d.Content.Find.Execute("22-3=", $true, $true, $false, $false, $false, $true, 1, $false, "27+65=", 2) | Out-Null
$d.Content.Find.Execute("36-29=", $true, $true, $false, $false, $false, $true, 1, $false, "50+30=", 2) | Out-Null
$d.Content.Find.Execute("39-29=", $true, $true, $false, $false, $false, $true, 1, $false, "53+37=", 2) | Out-Null
$d.Content.Find.Execute("14+23=", $true, $true, $false, $false, $false, $true, 1, $false, "57-7=", 2) | Out-Null
$d.Content.Find.Execute("78-58=", $true, $true, $false, $false, $false, $true, 1, $false, "67-42=", 2) | Out-Null
$d.Content.Find.Execute("3+4=", $true, $true, $false, $false, $false, $true, 1, $false, "66-58=", 2) | Out-Null
$d.Content.Find.Execute("42-16=", $true, $true, $false, $false, $false, $true, 1, $false, "74-0=", 2) | Out-Null
$d.Content.Find.Execute("17+12=", $true, $true, $false, $false, $false, $true, 1, $false, "80-23=", 2) | Out-Null
$d.Content.Find.Execute("75-69=", $true, $true, $false, $false, $false, $true, 1, $false, "28+4=", 2) | Out-Null
$d.Content.Find.Execute("3+37=", $true, $true, $false, $false, $false, $true, 1, $false, "2+77=", 2) | Out-Null
$d.Content.Find.Execute("32+19=", $true, $true, $false, $false, $false, $true, 1, $false, "37-13=", 2) | Out-Null
$d.Content.Find.Execute("18-15=", $true, $true, $false, $false, $false, $true, 1, $false, "10+19=", 2) | Out-Null
$d.Content.Find.Execute("45-30=", $true, $true, $false, $false, $false, $true, 1, $false, "28+23=", 2) | Out-Null
$d.Content.Find.Execute("60+11=", $true, $true, $false, $false, $false, $true, 1, $false, "84-32=", 2) | Out-Null
$d.Content.Find.Execute("9+12=", $true, $true, $false, $false, $false, $true, 1, $false, "66-36=", 2) | Out-Null
$d.Content.Find.Execute("39+38=", $true, $true, $false, $false, $false, $true, 1, $false, "13+50=", 2) | Out-Null
$d.Content.Find.Execute("19+20=", $true, $true, $false, $false, $false, $true, 1, $false, "35-10=", 2) | Out-Null
$d.Content.Find.Execute("13+63=", $true, $true, $false, $false, $false, $true, 1, $false, "39+17=", 2) | Out-Null
$d.Content.Find.Execute("23+25=", $true, $true, $false, $false, $false, $true, 1, $false, "46+1=", 2) | Out-Null
$d.Content.Find.Execute("0+80=", $true, $true, $false, $false, $false, $true, 1, $false, "71-31=", 2) | Out-Null
$d.Content.Find.Execute("33-10=", $true, $true, $false, $false, $false, $true, 1, $false, "51+8=", 2) | Out-Null
$d.Content.Find.Execute("11+81=", $true, $true, $false, $false, $false, $true, 1, $false, "84-71=", 2) | Out-Null
$d.Content.Find.Execute("10+10=", $true, $true, $false, $false, $false, $true, 1, $false, "98-2=", 2) | Out-Null
$d.Content.Find.Execute("41+35=", $true, $true, $false, $false, $false, $true, 1, $false, "78+3=", 2) | Out-Null
$d.Content.Find.Execute("34-19=", $true, $true, $false, $false, $false, $true, 1, $false, "31-25=", 2) | Out-Null
$d.Content.Find.Execute("92-6=", $true, $true, $false, $false, $false, $true, 1, $false, "47+3=", 2) | Out-Null
$d.Content.Find.Execute("3+17=", $true, $true, $false, $false, $false, $true, 1, $false, "27+47=", 2) | Out-Null
$d.Content.Find.Execute("59-52=", $true, $true, $false, $false, $false, $true, 1, $false, "56-4=", 2) | Out-Null
$d.Content.Find.Execute("82-50=", $true, $true, $false, $false, $false, $true, 1, $false, "12-0=", 2) | Out-Null
$d.Content.Find.Execute("11+10=", $true, $true, $false, $false, $false, $true, 1, $false, "95-27=", 2) | Out-Null
$d.Content.Find.Execute("98-21=", $true, $true, $false, $false, $false, $true, 1, $false, "41-11=", 2) | Out-Null
$d.Content.Find.Execute("71-29=", $true, $true, $false, $false, $false, $true, 1, $false, "91-42=", 2) | Out-Null
$d.Content.Find.Execute("82+3=", $true, $true, $false, $false, $false, $true, 1, $false, "46+40=", 2) | Out-Null
$d.Content.Find.Execute("70-55=", $true, $true, $false, $false, $false, $true, 1, $false, "48+33=", 2) | Out-Null
$d.Content.Find.Execute("72-40=", $true, $true, $false, $false, $false, $true, 1, $false, "50+45=", 2) | Out-Null
$d.Content.Find.Execute("11+65=", $true, $true, $false, $false, $false, $true, 1, $false, "89-34=", 2) | Out-Null
$d.Content.Find.Execute("81+2=", $true, $true, $false, $false, $false, $true, 1, $false, "55-40=", 2) | Out-Null
$d.Content.Find.Execute("58-40=", $true, $true, $false, $false, $false, $true, 1, $false, "25-1=", 2) | Out-Null
$d.Content.Find.Execute("85-35=", $true, $true, $false, $false, $false, $true, 1, $false, "90-55=", 2) | Out-Null
$d.Content.Find.Execute("37+6=", $true, $true, $false, $false, $false, $true, 1, $false, "13+6=", 2) | Out-Null
$d.Content.Find.Execute("50-40=", $true, $true, $false, $false, $false, $true, 1, $false, "33+48=", 2) | Out-Null
$d.Content.Find.Execute("43+36=", $true, $true, $false, $false, $false, $true, 1, $false, "27+53=", 2) | Out-Null
$d.Content.Find.Execute("49+17=", $true, $true, $false, $false, $false, $true, 1, $false, "64+33=", 2) | Out-Null
$d.Content.Find.Execute("84-16=", $true, $true, $false, $false, $false, $true, 1, $false, "84-7=", 2) | Out-Null
$d.Content.Find.Execute("55-29=", $true, $true, $false, $false, $false, $true, 1, $false, "88+0=", 2) | Out-Null
$d.Content.Find.Execute("48+36=", $true, $true, $false, $false, $false, $true, 1, $false, "35-30=", 2) | Out-Null
$d.Content.Find.Execute("58+21=", $true, $true, $false, $false, $false, $true, 1, $false, "87-2=", 2) | Out-Null
$d.Content.Find.Execute("51+28=", $true, $true, $false, $false, $false, $true, 1, $false, "65-18=", 2) | Out-Null
$d.Content.Find.Execute("9+64=", $true, $true, $false, $false, $false, $true, 1, $false, "13+68=", 2) | Out-Null
$d.Content.Find.Execute("69+24=", $true, $true, $false, $false, $false, $true, 1, $false, "46-36=", 2) | Out-Null
$d.Content.Find.Execute("26+50=", $true, $true, $false, $false, $false, $true, 1, $false, "48-34=", 2) | Out-Null
$d.Content.Find.Execute("99-23=", $true, $true, $false, $false, $false, $true, 1, $false, "66-51=", 2) | Out-Null
$d.Content.Find.Execute("70-43=", $true, $true, $false, $false, $false, $true, 1, $false, "13+50=", 2) | Out-Null
